$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-23 from 45170 to 45174
for ($row = 2; $row -le 23; $row++) {
    $ws.Cells.Item($row, 3).Value = 45174
}
